$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-record data which gets re-shuffled across rows 7-13:
# A(1) B(2) D(4) E(5) F(6) G(7) H(8) P(16) Q(17) R(18) AC(29)
$cols = @(1, 2, 4, 5, 6, 7, 8, 16, 17, 18, 29)

# Snapshot the current ("before") values for every affected row/column first,
# since the new row order reads from several different source rows and this
# is a cyclic permutation of rows (we must not overwrite a source row before
# it has been read). Read with Value2 (plain data, no formatting wrapper).
$snapshot = @{}
foreach ($r in 7..13) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# New row r takes its data from old row $map[r].
$map = @{
    7  = 12
    8  = 7
    9  = 13
    10 = 8
    11 = 9
    12 = 10
    13 = 11
}

foreach ($r in 7..13) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $val = $snapshot["$src,$c"]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).Value = ""
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
